$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-24 12:45:28"

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
